$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.296
$ws.Range("C4").Value = 0.051
$ws.Range("E4").Value = 0.194
$ws.Range("F4").Value = 0.027
$ws.Range("G4").Value = 0.166
$ws.Range("H4").Value = 0.214
$ws.Range("J4").Value = 0.157
$ws.Range("K4").Value = 0.346
$ws.Range("L4").Value = 0.098
$ws.Range("M4").Value = 0.313
$ws.Range("N4").Value = 0.271
$ws.Range("O4").Value = 0.02
$ws.Range("P4").Value = 0.143
$ws.Range("Z4").Value = 0.462
$ws.Range("AA4").Value = 0.125
$ws.Range("AB4").Value = 0.354
$ws.Range("AC4").Value = 0.12
$ws.Range("AE4").Value = 0.081
$ws.Range("AF4").Value = 0.737
$ws.Range("AG4").Value = 0.094
$ws.Range("AH4").Value = 0.307
$ws.Range("AI4").Value = 0.658
$ws.Range("AJ4").Value = 0.172
$ws.Range("AK4").Value = 0.415
$ws.Range("AL4").Value = 0.703
$ws.Range("AN4").Value = 0.339
$ws.Range("AO4").Value = 0.699
$ws.Range("B5").Value = 0.8159999999999999
$ws.Range("C5").Value = 0.15
$ws.Range("D5").Value = 0.388
$ws.Range("E5").Value = 0.769
$ws.Range("F5").Value = 0.178
$ws.Range("G5").Value = 0.421
$ws.Range("H5").Value = 0.857
$ws.Range("I5").Value = 0.122
$ws.Range("J5").Value = 0.35
$ws.Range("K5").Value = 0.675
$ws.Range("L5").Value = 0.219
$ws.Range("M5").Value = 0.468
$ws.Range("N5").Value = 0.842
$ws.Range("O5").Value = 0.133
$ws.Range("P5").Value = 0.365
$ws.Range("Z5").Value = 0.857
$ws.Range("AA5").Value = 0.122
$ws.Range("AB5").Value = 0.35
$ws.Range("AC5").Value = 0.727
$ws.Range("AD5").Value = 0.198
$ws.Range("AE5").Value = 0.445
$ws.Range("AF5").Value = 0.974
$ws.Range("AH5").Value = 0.16
$ws.Range("AI5").Value = 0.763
$ws.Range("AJ5").Value = 0.181
$ws.Range("AK5").Value = 0.425
$ws.Range("AL5").Value = 0.921
$ws.Range("AM5").Value = 0.073
$ws.Range("AN5").Value = 0.27
$ws.Range("AO5").Value = 0.886
$ws.Range("B6").Value = 0.434
$ws.Range("E6").Value = 0.31
$ws.Range("H6").Value = 0.342
$ws.Range("K6").Value = 0.457
$ws.Range("N6").Value = 0.41
$ws.Range("Z6").Value = 0.6
$ws.Range("AC6").Value = 0.206
$ws.Range("AF6").Value = 0.839
$ws.Range("AI6").Value = 0.707
$ws.Range("AL6").Value = 0.797
$ws.Range("AO6").Value = 0.781
$ws.Range("B7").Value = 0.604
$ws.Range("E7").Value = 0.483
$ws.Range("H7").Value = 0.535
$ws.Range("K7").Value = 0.5669999999999999
$ws.Range("N7").Value = 0.592
$ws.Range("Z7").Value = 0.732
$ws.Range("AC7").Value = 0.361
$ws.Range("AF7").Value = 0.915
$ws.Range("AI7").Value = 0.739
$ws.Range("AL7").Value = 0.867
$ws.Range("AO7").Value = 0.84
$ws.Range("B8").Value = 0.751
$ws.Range("C8").Value = 0.15
$ws.Range("D8").Value = 0.387
$ws.Range("E8").Value = 0.649
$ws.Range("F8").Value = 0.166
$ws.Range("G8").Value = 0.407
$ws.Range("H8").Value = 0.747
$ws.Range("I8").Value = 0.132
$ws.Range("J8").Value = 0.364
$ws.Range("K8").Value = 0.597
$ws.Range("L8").Value = 0.198
$ws.Range("M8").Value = 0.445
$ws.Range("N8").Value = 0.749
$ws.Range("O8").Value = 0.137
$ws.Range("P8").Value = 0.371
$ws.Range("Z8").Value = 0.789
$ws.Range("AA8").Value = 0.127
$ws.Range("AB8").Value = 0.356
$ws.Range("AC8").Value = 0.609
$ws.Range("AD8").Value = 0.185
$ws.Range("AE8").Value = 0.43
$ws.Range("AF8").Value = 0.893
$ws.Range("AG8").Value = 0.046
$ws.Range("AH8").Value = 0.215
$ws.Range("AI8").Value = 0.753
$ws.Range("AJ8").Value = 0.18
$ws.Range("AK8").Value = 0.424
$ws.Range("AL8").Value = 0.892
$ws.Range("AM8").Value = 0.078
$ws.Range("AN8").Value = 0.279
$ws.Range("AO8").Value = 0.846
$ws.Range("B9").Value = 0.658
$ws.Range("C9").Value = 0.225
$ws.Range("D9").Value = 0.474
$ws.Range("E9").Value = 0.513
$ws.Range("H9").Value = 0.619
$ws.Range("I9").Value = 0.236
$ws.Range("J9").Value = 0.486
$ws.Range("K9").Value = 0.5
$ws.Range("N9").Value = 0.632
$ws.Range("O9").Value = 0.233
$ws.Range("P9").Value = 0.482
$ws.Range("Z9").Value = 0.6899999999999999
$ws.Range("AA9").Value = 0.214
$ws.Range("AB9").Value = 0.462
$ws.Range("AC9").Value = 0.5
$ws.Range("AF9").Value = 0.763
$ws.Range("AG9").Value = 0.181
$ws.Range("AH9").Value = 0.425
$ws.Range("AI9").Value = 0.737
$ws.Range("AJ9").Value = 0.194
$ws.Range("AK9").Value = 0.44
$ws.Range("AL9").Value = 0.842
$ws.Range("AM9").Value = 0.133
$ws.Range("AN9").Value = 0.365
$ws.Range("AO9").Value = 0.781
$ws.Range("B10").Value = 0.8159999999999999
$ws.Range("C10").Value = 0.15
$ws.Range("D10").Value = 0.388
$ws.Range("E10").Value = 0.6919999999999999
$ws.Range("F10").Value = 0.213
$ws.Range("G10").Value = 0.462
$ws.Range("H10").Value = 0.786
$ws.Range("I10").Value = 0.168
$ws.Range("J10").Value = 0.41
$ws.Range("K10").Value = 0.675
$ws.Range("L10").Value = 0.219
$ws.Range("M10").Value = 0.468
$ws.Range("N10").Value = 0.8159999999999999
$ws.Range("O10").Value = 0.15
$ws.Range("P10").Value = 0.388
$ws.Range("Z10").Value = 0.857
$ws.Range("AA10").Value = 0.122
$ws.Range("AB10").Value = 0.35
$ws.Range("AC10").Value = 0.614
$ws.Range("AD10").Value = 0.237
$ws.Range("AE10").Value = 0.487
$ws.Range("AF10").Value = 0.974
$ws.Range("AH10").Value = 0.16
$ws.Range("AI10").Value = 0.763
$ws.Range("AJ10").Value = 0.181
$ws.Range("AK10").Value = 0.425
$ws.Range("AL10").Value = 0.921
$ws.Range("AM10").Value = 0.073
$ws.Range("AN10").Value = 0.27
$ws.Range("AO10").Value = 0.886
$ws.Range("B11").Value = 0.8159999999999999
$ws.Range("C11").Value = 0.15
$ws.Range("D11").Value = 0.388
$ws.Range("E11").Value = 0.769
$ws.Range("F11").Value = 0.178
$ws.Range("G11").Value = 0.421
$ws.Range("H11").Value = 0.857
$ws.Range("I11").Value = 0.122
$ws.Range("J11").Value = 0.35
$ws.Range("K11").Value = 0.675
$ws.Range("L11").Value = 0.219
$ws.Range("M11").Value = 0.468
$ws.Range("N11").Value = 0.842
$ws.Range("O11").Value = 0.133
$ws.Range("P11").Value = 0.365
$ws.Range("Z11").Value = 0.857
$ws.Range("AA11").Value = 0.122
$ws.Range("AB11").Value = 0.35
$ws.Range("AC11").Value = 0.659
$ws.Range("AD11").Value = 0.225
$ws.Range("AE11").Value = 0.474
$ws.Range("AF11").Value = 0.974
$ws.Range("AH11").Value = 0.16
$ws.Range("AI11").Value = 0.763
$ws.Range("AJ11").Value = 0.181
$ws.Range("AK11").Value = 0.425
$ws.Range("AL11").Value = 0.921
$ws.Range("AM11").Value = 0.073
$ws.Range("AN11").Value = 0.27
$ws.Range("AO11").Value = 0.886
$ws.Range("B12").Value = 1.258
$ws.Range("C12").Value = 0.32
$ws.Range("D12").Value = 0.5659999999999999
$ws.Range("E12").Value = 1.633
$ws.Range("F12").Value = 1.032
$ws.Range("G12").Value = 1.016
$ws.Range("H12").Value = 1.556
$ws.Range("I12").Value = 1.191
$ws.Range("J12").Value = 1.091
$ws.Range("K12").Value = 1.407
$ws.Range("L12").Value = 0.538
$ws.Range("M12").Value = 0.733
$ws.Range("N12").Value = 1.406
$ws.Range("O12").Value = 0.616
$ws.Range("P12").Value = 0.785
$ws.Range("Z12").Value = 1.25
$ws.Range("AA12").Value = 0.299
$ws.Range("AB12").Value = 0.546
$ws.Range("AC12").Value = 2
$ws.Range("AD12").Value = 3.812
$ws.Range("AE12").Value = 1.953
$ws.Range("AF12").Value = 1.243
$ws.Range("AG12").Value = 0.238
$ws.Range("AH12").Value = 0.488
$ws.Range("AI12").Value = 1.034
$ws.Range("AJ12").Value = 0.033
$ws.Range("AK12").Value = 0.182
$ws.Range("AL12").Value = 1.086
$ws.Range("AM12").Value = 0.078
$ws.Range("AN12").Value = 0.28
$ws.Range("AO12").Value = 1.121
$ws.Range("B13").Value = 3.474
$ws.Range("C13").Value = 1.46
$ws.Range("D13").Value = 1.208
$ws.Range("E13").Value = 4.564
$ws.Range("F13").Value = 0.707
$ws.Range("G13").Value = 0.841
$ws.Range("H13").Value = 4.524
$ws.Range("I13").Value = 0.916
$ws.Range("J13").Value = 0.957
$ws.Range("K13").Value = 2.3
$ws.Range("L13").Value = 0.61
$ws.Range("M13").Value = 0.781
$ws.Range("N13").Value = 3.263
$ws.Range("O13").Value = 0.72
$ws.Range("P13").Value = 0.849
$ws.Range("Z13").Value = 2.833
$ws.Range("AA13").Value = 3.901
$ws.Range("AB13").Value = 1.975
$ws.Range("AC13").Value = 6.273
$ws.Range("AD13").Value = 2.88
$ws.Range("AE13").Value = 1.697
$ws.Range("AF13").Value = 1.605
$ws.Range("AG13").Value = 0.713
$ws.Range("AH13").Value = 0.844
$ws.Range("AI13").Value = 1.289
$ws.Range("AJ13").Value = 0.364
$ws.Range("AK13").Value = 0.603
$ws.Range("AL13").Value = 1.579
$ws.Range("AM13").Value = 0.717
$ws.Range("AN13").Value = 0.847
$ws.Range("AO13").Value = 1.491
